$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 157 (pushes old rows 157:163 down to 160:166)
$ws.Rows("157:159").Insert()

# New row 157: Damasco / Dina / Especial - San Felipe de Aconcagua, week of 2023-01-05
$ws.Cells.Item(157, 1).Value = 9
$ws.Cells.Item(157, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(157, 3).Value = "Metropolitana"
$ws.Cells.Item(157, 4).Value = 44931
$ws.Cells.Item(157, 5).Value = 13
$ws.Cells.Item(157, 6).Value = "Fruta"
$ws.Cells.Item(157, 7).Value = 100103
$ws.Cells.Item(157, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(157, 9).Value = 100103003
$ws.Cells.Item(157, 10).Value = "Damasco"
$ws.Cells.Item(157, 11).Value = "Dina"
$ws.Cells.Item(157, 12).Value = "Especial"
$ws.Cells.Item(157, 13).Value = 330
$ws.Cells.Item(157, 14).Value = 16000
$ws.Cells.Item(157, 15).Value = 16000
$ws.Cells.Item(157, 16).Value = 16000
$ws.Cells.Item(157, 17).Value = "$/caja 16 kilos granel"
$ws.Cells.Item(157, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(157, 19).Value = 1000
$ws.Cells.Item(157, 20).Value = 16

# New row 158: Damasco / Dina / Primera - San Felipe de Aconcagua, week of 2023-01-05
$ws.Cells.Item(158, 1).Value = 9
$ws.Cells.Item(158, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(158, 3).Value = "Metropolitana"
$ws.Cells.Item(158, 4).Value = 44931
$ws.Cells.Item(158, 5).Value = 13
$ws.Cells.Item(158, 6).Value = "Fruta"
$ws.Cells.Item(158, 7).Value = 100103
$ws.Cells.Item(158, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(158, 9).Value = 100103003
$ws.Cells.Item(158, 10).Value = "Damasco"
$ws.Cells.Item(158, 11).Value = "Dina"
$ws.Cells.Item(158, 12).Value = "Primera"
$ws.Cells.Item(158, 13).Value = 310
$ws.Cells.Item(158, 14).Value = 12000
$ws.Cells.Item(158, 15).Value = 12000
$ws.Cells.Item(158, 16).Value = 12000
$ws.Cells.Item(158, 17).Value = "$/caja 16 kilos granel"
$ws.Cells.Item(158, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(158, 19).Value = 750
$ws.Cells.Item(158, 20).Value = 16

# New row 159: Damasco / Dina / Segunda - San Felipe de Aconcagua, week of 2023-01-05
$ws.Cells.Item(159, 1).Value = 9
$ws.Cells.Item(159, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(159, 3).Value = "Metropolitana"
$ws.Cells.Item(159, 4).Value = 44931
$ws.Cells.Item(159, 5).Value = 13
$ws.Cells.Item(159, 6).Value = "Fruta"
$ws.Cells.Item(159, 7).Value = 100103
$ws.Cells.Item(159, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(159, 9).Value = 100103003
$ws.Cells.Item(159, 10).Value = "Damasco"
$ws.Cells.Item(159, 11).Value = "Dina"
$ws.Cells.Item(159, 12).Value = "Segunda"
$ws.Cells.Item(159, 13).Value = 280
$ws.Cells.Item(159, 14).Value = 9600
$ws.Cells.Item(159, 15).Value = 9600
$ws.Cells.Item(159, 16).Value = 9600
$ws.Cells.Item(159, 17).Value = "$/caja 16 kilos granel"
$ws.Cells.Item(159, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(159, 19).Value = 600
$ws.Cells.Item(159, 20).Value = 16
